$d = $word.ActiveDocument

function Get-ExactRange([string]$needle, [int]$startAt) {
    # Locate an exact, case-sensitive match of $needle starting the search
    # at character position $startAt. Returns a fresh Range over the match
    # (so callers can safely mutate it with InsertXML).
    $docEnd = $d.Content.End
    $searchRange = $d.Range($startAt, $docEnd)
    $ok = $searchRange.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $needle"
    }
    $foundStart = $searchRange.Start
    $foundEnd = $searchRange.End
    return $d.Range($foundStart, $foundEnd)
}

function Insert-ParagraphsXml([int]$start, [int]$end, [string]$bodyFragment) {
    # Replace the (start,end) range -- which should span whole paragraphs,
    # paragraph mark included -- with the <w:p>...</w:p> fragment(s)
    # supplied in $bodyFragment.
    $rng = $d.Range($start, $end)
    $header = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
    $footer = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $xml = $header + $bodyFragment + $footer
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1) "Attributes of a collection:" -> "CollectionID" (flagged as a
#    misspelling by the proofer, hence the spellStart/spellEnd wrap).
# ---------------------------------------------------------------------
$r1 = Get-ExactRange "Attributes of a collection:" 0
$frag1 = '<w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>CollectionID</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$start1 = $r1.Start
$end1 = $r1.End
Insert-ParagraphsXml $start1 $end1 $frag1

# ---------------------------------------------------------------------
# 2) The two "Collection" attribute bullets covering influencer/IP
#    collaborations are reworded and doubled (Is.../...Name pairs).
# ---------------------------------------------------------------------
$rStart2 = Get-ExactRange "Influencer Collaboration? (Yes/No)" 0
$rEnd2 = Get-ExactRange "Influencer_Name" 0
$start2 = $rStart2.Start
# Extend through the paragraph mark that ends the "Influencer_Name" paragraph.
$end2 = $rEnd2.End + 1
$frag2 = (
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>IsInfluencerCollaboration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (Yes/No)</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>InfluencerCollaborationName</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (NULL)</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>IsIPCollaboration</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (Yes/No)</w:t></w:r></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>IPCollaborationName</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (NULL)</w:t></w:r></w:p>'
)
Insert-ParagraphsXml $start2 $end2 $frag2

# ---------------------------------------------------------------------
# 3) PolishFinish: "DryTexture (glossy, satin, matte)" -> split into
#    "IsDriesMatte" and a brand new "IsGlowitheDark" bullet.
# ---------------------------------------------------------------------
$rStart3 = Get-ExactRange "DryTexture" 0
$rEnd3 = Get-ExactRange " (glossy, satin, matte)" 0
$start3 = $rStart3.Start
$end3 = $rEnd3.End + 1
$frag3 = (
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>IsDriesMatte</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>' +
  '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>IsGlowitheDark</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
)
Insert-ParagraphsXml $start3 $end3 $frag3

# ---------------------------------------------------------------------
# 4) The rendered-page-break marker shifts from the "Organizer_ID" run
#    up to the "Theme" run (earlier insertions pushed the page boundary
#    up by one run).
# ---------------------------------------------------------------------
# Find the second "Theme" occurrence -- the one that immediately precedes
# "Organizer_id (Organizing group)".
$themeSearchStart = 0
$found4 = $false
while (-not $found4) {
    $rTheme = Get-ExactRange "Theme" $themeSearchStart
    $afterStart = $rTheme.End
    $afterEndCandidate = $afterStart + 40
    $docEnd4 = $d.Content.End
    $afterEnd = [Math]::Min($afterEndCandidate, $docEnd4)
    $rAfter = $d.Range($afterStart, $afterEnd)
    $afterText = $rAfter.Text
    if ($afterText -like "*Organizer_id*") {
        $found4 = $true
    } else {
        $themeSearchStart = $rTheme.End
    }
}

$themeStart = $rTheme.Start
$themeEnd = $rTheme.End + 1
$themeFrag = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Theme</w:t></w:r></w:p>'
Insert-ParagraphsXml $themeStart $themeEnd $themeFrag

$rOrganizer2 = Get-ExactRange "Organizer_ID" 0
$organizerStart = $rOrganizer2.Start
$organizerEnd = $rOrganizer2.End + 1
$organizerFrag = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Organizer_ID</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
Insert-ParagraphsXml $organizerStart $organizerEnd $organizerFrag

Write-Host "Edits applied"
